$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 0.7889997959136963
$ws.Cells.Item(3, 1).Value = 0.13100004196166992
$ws.Cells.Item(4, 1).Value = 0.15399980545043945
$ws.Cells.Item(5, 1).Value = 0.09200000762939453
$ws.Cells.Item(6, 1).Value = 0.0839998722076416
$ws.Cells.Item(7, 1).Value = 0.0840001106262207
$ws.Cells.Item(8, 1).Value = 0.09800004959106445
$ws.Cells.Item(9, 1).Value = 0.11100006103515625
$ws.Cells.Item(10, 1).Value = 0.09100008010864258
$ws.Cells.Item(11, 1).Value = 0.09500002861022949
$ws.Cells.Item(12, 1).Value = 0.08699989318847656
$ws.Cells.Item(13, 1).Value = 0.09099984169006348
$ws.Cells.Item(14, 1).Value = 0.1119999885559082
$ws.Cells.Item(15, 1).Value = 0.10899996757507324
$ws.Cells.Item(16, 1).Value = 0.10100007057189941
$ws.Cells.Item(17, 1).Value = 0.09400010108947754
$ws.Cells.Item(18, 1).Value = 0.0820000171661377
$ws.Cells.Item(19, 1).Value = 0.08999991416931152
$ws.Cells.Item(20, 1).Value = 0.08499979972839355
$ws.Cells.Item(21, 1).Value = 0.12400007247924805
$ws.Cells.Item(22, 1).Value = 0.1119999885559082
$ws.Cells.Item(23, 1).Value = 0.09599995613098145
$ws.Cells.Item(24, 1).Value = 0.09300017356872559
$ws.Cells.Item(25, 1).Value = 0.08300018310546875
$ws.Cells.Item(26, 1).Value = 0.09599995613098145
$ws.Cells.Item(27, 1).Value = 0.09200000762939453
$ws.Cells.Item(28, 1).Value = 0.10299992561340332
$ws.Cells.Item(29, 1).Value = 0.15899991989135742
$ws.Cells.Item(30, 1).Value = 0.11800003051757812
$ws.Cells.Item(31, 1).Value = 0.1119999885559082
$ws.Cells.Item(32, 1).Value = 0.09599995613098145
$ws.Cells.Item(33, 1).Value = 0.08100008964538574
$ws.Cells.Item(34, 1).Value = 0.0840001106262207
$ws.Cells.Item(35, 1).Value = 0.0989999771118164
$ws.Cells.Item(36, 1).Value = 0.08699989318847656
$ws.Cells.Item(37, 1).Value = 0.08299994468688965
$ws.Cells.Item(38, 1).Value = 0.09599995613098145
$ws.Cells.Item(39, 1).Value = 0.09299993515014648
$ws.Cells.Item(40, 1).Value = 0.11299991607666016
$ws.Cells.Item(41, 1).Value = 0.08299994468688965
$ws.Cells.Item(42, 1).Value = 0.11599993705749512
$ws.Cells.Item(43, 1).Value = 0.09400010108947754
$ws.Cells.Item(44, 1).Value = 0.10100007057189941
$ws.Cells.Item(45, 1).Value = 0.08799982070922852
$ws.Cells.Item(46, 1).Value = 0.11299991607666016
$ws.Cells.Item(47, 1).Value = 0.09200000762939453
$ws.Cells.Item(48, 1).Value = 0.10300016403198242
$ws.Cells.Item(49, 1).Value = 0.08699989318847656
$ws.Cells.Item(50, 1).Value = 0.08100008964538574
$ws.Cells.Item(51, 1).Value = 0.11399984359741211
$ws.Cells.Item(52, 1).Value = 0.09200000762939453
$ws.Cells.Item(53, 1).Value = 0.09600019454956055
$ws.Cells.Item(54, 1).Value = 0.09100008010864258
$ws.Cells.Item(55, 1).Value = 0.08800005912780762
$ws.Cells.Item(56, 1).Value = 0.0970001220703125
$ws.Cells.Item(57, 1).Value = 0.11100006103515625
$ws.Cells.Item(58, 1).Value = 0.10000014305114746
$ws.Cells.Item(59, 1).Value = 0.09100008010864258
$ws.Cells.Item(60, 1).Value = 0.0970001220703125
$ws.Cells.Item(61, 1).Value = 0.08500003814697266
$ws.Cells.Item(62, 1).Value = 0.09399986267089844
$ws.Cells.Item(63, 1).Value = 0.10800004005432129
$ws.Cells.Item(64, 1).Value = 0.10699987411499023
$ws.Cells.Item(65, 1).Value = 0.08700013160705566
$ws.Cells.Item(66, 1).Value = 0.08599996566772461
$ws.Cells.Item(67, 1).Value = 0.0820000171661377
$ws.Cells.Item(68, 1).Value = 0.09400010108947754
$ws.Cells.Item(69, 1).Value = 0.1099998950958252
$ws.Cells.Item(70, 1).Value = 0.0839998722076416
$ws.Cells.Item(71, 1).Value = 0.10199999809265137
$ws.Cells.Item(72, 1).Value = 0.0989999771118164
$ws.Cells.Item(73, 1).Value = 0.09099984169006348
$ws.Cells.Item(74, 1).Value = 0.09400010108947754
$ws.Cells.Item(75, 1).Value = 0.08999991416931152
$ws.Cells.Item(76, 1).Value = 0.10500001907348633
$ws.Cells.Item(77, 1).Value = 0.08600020408630371
$ws.Cells.Item(78, 1).Value = 0.07999992370605469
$ws.Cells.Item(79, 1).Value = 0.07899999618530273
$ws.Cells.Item(80, 1).Value = 0.0969998836517334
$ws.Cells.Item(81, 1).Value = 0.09800004959106445
